# Update the "取得日時" (retrieved-at) timestamp in column A for all data
# rows (2-15) of the "ランサーズ" sheet to the new run timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-23 18:39:47"

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
